$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.856.57'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.794.47'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.74'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.65'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.552'
$ws.Range("E7").Value = '  -2.61%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.598'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.05'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("E11").Value = '  +2.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.23'
$ws.Range("E12").Value = '  +3.41%  '
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.66'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '3.233.94'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '2.827.03'
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '51.820.64'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.74'
$ws.Range("E19").Value = '  +4.11%  '
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.19'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("D22").Value = '0.0₃0967'
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.12'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.00'
$ws.Range("E24").Value = '  -2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.74'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.24'
$ws.Range("E26").Value = '  -1.86%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.163'
$ws.Range("E28").Value = '  +11.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.26'
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.10'
$ws.Range("E30").Value = '  +8.63%  '
$ws.Range("E31").Value = '  +9.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.01'
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0456'
$ws.Range("E33").Value = '  -2.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.59'
$ws.Range("E34").Value = '  +5.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.95'
$ws.Range("E35").Value = '  -12.07%  '
$ws.Range("E36").Value = '  -1.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.59'
$ws.Range("E38").Value = '  +2.68%  '
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  -1.35%  '
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.114'
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '121.23'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.14'
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").Value = '2.142.75'
$ws.Range("E46").Value = '  +3.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.32'
$ws.Range("E48").Value = '  +6.34%  '
$ws.Range("B49").Value = 'SEI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.914'
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.47'
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("E51").Value = '  +8.86%  '
